$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 27
$ws.Range("H2").Value = 27

# Row 7
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 3
$ws.Range("H7").Value = 3

# Row 9
$ws.Range("E9").Value = 25
$ws.Range("F9").Value = 9
$ws.Range("H9").Value = 9

# Row 16
$ws.Range("E16").Value = 11

# Row 19
$ws.Range("E19").Value = 54
$ws.Range("F19").Value = 27
$ws.Range("H19").Value = 27

# Row 25
$ws.Range("F25").Value = 10
$ws.Range("H25").Value = 10

# Row 27
$ws.Range("E27").Value = 9

# Row 37
$ws.Range("E37").Value = 50

# Row 38
$ws.Range("E38").Value = 70

# Row 46
$ws.Range("E46").Value = 27

# Row 47
$ws.Range("F47").Value = 35
$ws.Range("H47").Value = 35

# Row 63
$ws.Range("F63").Value = 12
$ws.Range("H63").Value = 12

# Row 70
$ws.Range("F70").Value = 20
$ws.Range("H70").Value = 20

# Row 80
$ws.Range("F80").Value = 9
$ws.Range("H80").Value = 9

# Row 87
$ws.Range("E87").Value = 14
